$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contact name in A2 (shorten "Ashraful Hussain Khan" -> "Ashraful")
$ws.Range("A2").Value = "Ashraful"

# Remove row 3 entirely (the "Hafiz" contact), so the same contact (Ashraful)
# can now receive multiple messages instead of sharing rows with other contacts.
$ws.Rows("3").Delete()

# Move the selection, matching the saved view state in the workbook.
$ws.Range("C3").Select()
